$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value while always preserving it as literal text
# (coinranking price/volume strings must never be re-interpreted as numbers,
# which would silently drop trailing zeros, use exponential notation, or
# introduce floating point rounding noise).
function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    if ($text -match "^[+-]?\d+(\.\d+)?$") {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $text
}

# Row 2
Set-TextValue 2 4 "65.866.44"
Set-TextValue 2 5 "  -0.47%  "

# Row 3
Set-TextValue 3 4 "2.675.94"

# Row 4
Set-TextValue 4 5 "  +0.01%  "

# Row 5
Set-TextValue 5 4 "600.95"
Set-TextValue 5 5 "  -1.32%  "

# Row 6
Set-TextValue 6 4 "157.81"
Set-TextValue 6 5 "  -0.84%  "

# Row 7
Set-TextValue 7 5 "  +0.03%  "

# Row 8
Set-TextValue 8 5 "  +3.08%  "

# Row 9
Set-TextValue 9 4 "0.130"
Set-TextValue 9 5 "  +2.38%  "

# Row 10
Set-TextValue 10 4 "0.399"
Set-TextValue 10 5 "  -1.36%  "

# Row 11
Set-TextValue 11 4 "5.84"
Set-TextValue 11 5 "  -3.37%  "

# Row 12
Set-TextValue 12 5 "  -0.11%  "

# Row 13
Set-TextValue 13 4 "29.20"
Set-TextValue 13 5 "  -4.76%  "

# Row 14
Set-TextValue 14 4 "0.0000201"
Set-TextValue 14 5 "  -6.28%  "

# Row 15
Set-TextValue 15 4 "3.156.55"
Set-TextValue 15 5 "  -1.04%  "

# Row 16
Set-TextValue 16 4 "65.708.99"
Set-TextValue 16 5 "  -0.46%  "

# Row 17
Set-TextValue 17 4 "2.682.05"
Set-TextValue 17 5 "  -1.26%  "

# Row 18
Set-TextValue 18 5 "  -0.44%  "

# Row 19
Set-TextValue 19 4 "4.82"
Set-TextValue 19 5 "  -2.24%  "

# Row 20
Set-TextValue 20 4 "7.56"
Set-TextValue 20 5 "  -3.72%  "

# Row 21
Set-TextValue 21 4 "352.82"
Set-TextValue 21 5 "  -2.24%  "

# Row 22
Set-TextValue 22 4 "1.00"
Set-TextValue 22 5 "  -0.06%  "

# Row 23
Set-TextValue 23 4 "69.74"
Set-TextValue 23 5 "  -2.26%  "

# Row 24
Set-TextValue 24 4 "0.0000118"
Set-TextValue 24 5 "  +3.47%  "

# Row 25
Set-TextValue 25 4 "9.65"
Set-TextValue 25 5 "  -3.22%  "

# Row 26
Set-TextValue 26 5 "  +0.87%  "

# Row 27
Set-TextValue 27 4 "1.61"
Set-TextValue 27 5 "  -5.13%  "

# Row 28
Set-TextValue 28 4 "0.166"
Set-TextValue 28 5 "  -4.51%  "

# Row 29
Set-TextValue 29 4 "7.99"
Set-TextValue 29 5 "  -4.70%  "

# Row 30 now: PancakeSwap
Set-TextValue 30 2 "PancakeSwap"
Set-TextValue 30 3 "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue 30 4 "2.16"
Set-TextValue 30 5 "  -2.74%  "

# Row 31 now: Binance-PegBSC-USD
Set-TextValue 31 2 "Binance-PegBSC-USD"
Set-TextValue 31 3 "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue 31 4 "0.999"
Set-TextValue 31 5 "  -0.08%  "

# Row 32 now: Bittensor
Set-TextValue 32 2 "Bittensor"
Set-TextValue 32 3 "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue 32 4 "536.42"
Set-TextValue 32 5 "  -1.76%  "

# Row 33
Set-TextValue 33 4 "1.77"
Set-TextValue 33 5 "  -3.14%  "

# Row 34
Set-TextValue 34 4 "6.53"
Set-TextValue 34 5 "  -3.56%  "

# Row 35
Set-TextValue 35 4 "5.46"
Set-TextValue 35 5 "  -1.12%  "

# Row 36
Set-TextValue 36 4 "0.424"
Set-TextValue 36 5 "  -3.16%  "

# Row 37
Set-TextValue 37 4 "20.60"
Set-TextValue 37 5 "  -1.59%  "

# Row 38
Set-TextValue 38 5 "  -0.05%  "

# Row 39
Set-TextValue 39 4 "157.89"
Set-TextValue 39 5 "  -3.80%  "

# Row 40
Set-TextValue 40 4 "1.95"
Set-TextValue 40 5 "  -3.28%  "

# Row 41
Set-TextValue 41 4 "0.999"

# Row 42
Set-TextValue 42 4 "162.77"
Set-TextValue 42 5 "  -4.81%  "

# Row 43
Set-TextValue 43 4 "4.13"
Set-TextValue 43 5 "  -2.05%  "

# Row 44
Set-TextValue 44 4 "2.38"
Set-TextValue 44 5 "  +0.47%  "

# Row 45
Set-TextValue 45 4 "0.0610"
Set-TextValue 45 5 "  -2.82%  "

# Row 46
Set-TextValue 46 4 "22.73"
Set-TextValue 46 5 "  -4.88%  "

# Row 47
Set-TextValue 47 4 "0.642"
Set-TextValue 47 5 "  -3.01%  "

# Row 48
Set-TextValue 48 4 "0.0258"
Set-TextValue 48 5 "  -4.15%  "

# Row 49 now: BabyDogeCoin
Set-TextValue 49 2 "BabyDogeCoin"
Set-TextValue 49 3 "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue 49 4 "0.0₆0258"
Set-TextValue 49 5 "  +10.80%  "

# Row 50 now: EnergySwap
Set-TextValue 50 2 "EnergySwap"
Set-TextValue 50 3 "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue 50 4 "20.18"
Set-TextValue 50 5 "  -3.69%  "

# Row 51
Set-TextValue 51 4 "0.0991"
Set-TextValue 51 5 "  -0.46%  "
